{"js": "const body = context.document.body;\n\n// Locate the target paragraphs by their distinctive text.\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet pTable = null, pNote = null, pSenseTalk = null, pStrike = null, pRequire = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (pTable === null && t.indexOf(\"A table summarizing features\") === 0) {\n    pTable = paras.items[i];\n  } else if (pNote === null && t.indexOf(\"NOTE: Leaving table unformatted\") === 0) {\n    pNote = paras.items[i];\n  } else if (pSenseTalk === null && t.indexOf(\"default scripting language is based\") !== -1) {\n    pSenseTalk = paras.items[i];\n  } else if (pStrike === null && t.indexOf(\"scripts can easily be controlled with an source revision control system\") !== -1) {\n    pStrike = paras.items[i];\n  } else if (pRequire === null && t.indexOf(\"commercial license\") !== -1 && t.indexOf(\"Appendix A for details\") !== -1) {\n    pRequire = paras.items[i];\n  }\n}\n\n// 1) Turn the \"A table summarizing...\" paragraph red.\nif (pTable) {\n  pTable.font.color = \"#FF0000\";\n}\n\n// 2) Turn the \"NOTE: Leaving table unformatted...\" paragraph red (keeping its italics).\nif (pNote) {\n  pNote.font.color = \"#FF0000\";\n}\n\n// 3) \"manage tests\" -> \"manage test\"\nconst manageResults = body.search(\"manage tests\", { matchCase: true });\nmanageResults.load(\"items\");\nawait context.sync();\nif (manageResults.items.length > 0) {\n  manageResults.items[0].insertText(\"manage test\", Word.InsertLocation.replace);\n}\n\n// 4) Fix capitalization \"eggPlant's\" -> \"Eggplant's\" in the SenseTalk paragraph.\nif (pSenseTalk) {\n  const eggResults = pSenseTalk.search(\"eggPlant\", { matchCase: true });\n  eggResults.load(\"items\");\n  await context.sync();\n  if (eggResults.items.length > 0) {\n    eggResults.items[0].insertText(\"Eggplant\", Word.InsertLocation.replace);\n  }\n}\n\n// 5) Delete the whole strikethrough paragraph about source revision control.\nawait context.sync();\nif (pStrike) {\n  pStrike.delete();\n}\n\n// 6) \"eggPlant require\" -> \"Eggplant  require\" (note double space) plus a _GoBack bookmark\n//    right before \"require\", replicating the cursor position left behind by the edit.\nawait context.sync();\nif (pRequire) {\n  pRequire.load(\"text\");\n  await context.sync();\n  const reqResults = pRequire.search(\"eggPlant require\", { matchCase: true });\n  reqResults.load(\"items\");\n  await context.sync();\n  if (reqResults.items.length > 0) {\n    reqResults.items[0].insertText(\"Eggplant  require\", Word.InsertLocation.replace);\n  }\n}\n\n// 7) Move the _GoBack bookmark from the end of the \"Eggplant runs on a host computer...\"\n//    paragraph to just before \"require\" in the paragraph updated above.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nif (pRequire) {\n  const reqResults2 = pRequire.search(\"require\", { matchCase: true });\n  reqResults2.load(\"items\");\n  await context.sync();\n  if (reqResults2.items.length > 0) {\n    reqResults2.items[0].insertBookmark(\"_GoBack\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the target paragraphs by their distinctive text.\n$pTableIdx = -1\n$pNoteIdx = -1\n$pSenseTalkIdx = -1\n$pStrikeIdx = -1\n$pRequireIdx = -1\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($pTableIdx -eq -1 -and $t -like \"A table summarizing features*\") { $pTableIdx = $i }\n    elseif ($pNoteIdx -eq -1 -and $t -like \"NOTE: Leaving table unformatted*\") { $pNoteIdx = $i }\n    elseif ($pSenseTalkIdx -eq -1 -and $t -like \"*default scripting language is based*\") { $pSenseTalkIdx = $i }\n    elseif ($pStrikeIdx -eq -1 -and $t -like \"*scripts can easily be controlled with an source revision control system*\") { $pStrikeIdx = $i }\n    elseif ($pRequireIdx -eq -1 -and ($t -like \"*commercial license*\") -and ($t -like \"*Appendix A for details*\")) { $pRequireIdx = $i }\n}\n\n# 1) Turn the \"A table summarizing...\" paragraph red.\nif ($pTableIdx -ne -1) {\n    $d.Paragraphs.Item($pTableIdx).Range.Font.Color = 255\n}\n\n# 2) Turn the \"NOTE: Leaving table unformatted...\" paragraph red (keeping its italics).\nif ($pNoteIdx -ne -1) {\n    $d.Paragraphs.Item($pNoteIdx).Range.Font.Color = 255\n}\n\n# 3) \"manage tests\" -> \"manage test\"\n$rngManage = $d.Content\n$rngManage.Find.ClearFormatting()\n$rngManage.Find.Execute(\"manage tests\", $false, $true, $false, $false, $false, $true, 1, $false, \"manage test\", 2) | Out-Null\n\n# 4) Fix capitalization \"eggPlant's\" -> \"Eggplant's\" in the SenseTalk paragraph.\nif ($pSenseTalkIdx -ne -1) {\n    $rngSense = $d.Paragraphs.Item($pSenseTalkIdx).Range\n    $rngSense.Find.ClearFormatting()\n    $rngSense.Find.Execute(\"eggPlant\", $false, $true, $false, $false, $false, $true, 1, $false, \"Eggplant\", 2) | Out-Null\n}\n\n# 5) Delete the whole strikethrough paragraph about source revision control.\n#    (Re-resolve the index since earlier edits don't change paragraph counts before it.)\nif ($pStrikeIdx -ne -1) {\n    $d.Paragraphs.Item($pStrikeIdx).Range.Delete()\n}\n\n# 6) \"eggPlant require\" -> \"Eggplant  require\" (note double space).\nif ($pRequireIdx -ne -1) {\n    # The strike paragraph removal shifts indices for paragraphs after it by -1.\n    if ($pStrikeIdx -ne -1 -and $pStrikeIdx -lt $pRequireIdx) {\n        $pRequireIdx = $pRequireIdx - 1\n    }\n    $rngReq = $d.Paragraphs.Item($pRequireIdx).Range\n    $rngReq.Find.ClearFormatting()\n    $rngReq.Find.Execute(\"eggPlant require\", $false, $true, $false, $false, $false, $true, 1, $false, \"Eggplant  require\", 2) | Out-Null\n}\n\n# 7) Move the _GoBack bookmark from the end of the \"Eggplant runs on a host computer...\"\n#    paragraph to just before \"require\" in the paragraph updated above.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\nif ($pRequireIdx -ne -1) {\n    $rngFind = $d.Paragraphs.Item($pRequireIdx).Range.Duplicate()\n    $rngFind.Find.ClearFormatting()\n    $found = $rngFind.Find.Execute(\"require\")\n    if ($found) {\n        $rngFind.Collapse(1)\n        $d.Bookmarks.Add(\"_GoBack\", $rngFind) | Out-Null\n    }\n}\n"}
